$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the score cells in row 15 (student #12)
$ws.Range("D15").Value = 5
$ws.Range("E15").Value = 5
$ws.Range("F15").Value = 5
$ws.Range("G15").Value = 5

# Update the active selection / cell to I15, matching the saved sheet view
$ws.Range("I15").Select()
